$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- K3: header cell "2020" (bold Times New Roman 10, top+bottom medium border) ---
$ws.Range("K3").Value = 2020
$ws.Range("I3").Copy()
$ws.Range("K3").PasteSpecial(-4122)   # xlPasteFormats - reuse existing border (top+bottom medium)
$ws.Range("K3").Font.Name = "Times New Roman"
$ws.Range("K3").Font.Size = 10
$ws.Range("K3").Font.Bold = $true
$ws.Range("K3").VerticalAlignment = -4107   # xlVAlignBottom (default - clears inherited center)

# --- K4: data cell (Kyrghyz Times 9, right aligned, #,##0.0, top medium border only) ---
$ws.Range("K4").Value = 2.8218550629805335
$ws.Range("K4").Borders.Item(8).Weight = -4138   # xlEdgeTop, xlMedium
$ws.Range("K4").Font.Size = 9
$ws.Range("K4").Font.Name = "Kyrghyz Times"
$ws.Range("K4").NumberFormat = "#,##0.0"
$ws.Range("K4").HorizontalAlignment = -4152   # xlRight
$ws.Range("K4").VerticalAlignment = -4107     # xlVAlignBottom (default - clears inherited center)

# --- K5: data cell (Kyrghyz Times 9, right aligned, #,##0.0, bottom medium border) ---
$ws.Range("K5").Value = 1.3005071159823327
$ws.Range("A2").Copy()
$ws.Range("K5").PasteSpecial(-4122)   # xlPasteFormats - reuse existing border (bottom medium)
$ws.Range("K5").Font.Size = 9
$ws.Range("K5").Font.Name = "Kyrghyz Times"
$ws.Range("K5").NumberFormat = "#,##0.0"
$ws.Range("K5").HorizontalAlignment = -4152   # xlRight
$ws.Range("K5").VerticalAlignment = -4107     # xlVAlignBottom (default - clears inherited center)

# --- selection moves to L8 ---
$ws.Range("L8").Select()
